$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "2024-06-15 02:33:19"
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = 2

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "2024-06-15 02:33:19"
$ws.Range("D11").Value = 200
$ws.Range("E11").Value = 0
